$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1000
$ws.Range("I40").Value = 700
$ws.Range("K40").Value = 700
$ws.Range("M40").Value = -525

$ws.Range("H137").Value = 160485.36
$ws.Range("I137").Value = 3498.25
$ws.Range("J137").Value = 223280.2
$ws.Range("K137").Value = 10494.75
$ws.Range("L137").Value = 669840.6000000001
$ws.Range("M137").Value = -7944.75
$ws.Range("N137").Value = -674940.6000000001

$ws.Range("H138").Value = 2508.1428
$ws.Range("I138").Value = 1888.5555
$ws.Range("K138").Value = 5665.666499999999
$ws.Range("M138").Value = -525.6664999999994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 314.83334
$ws.Range("I5").Value = 293.8
$ws.Range("K5").Value = 293.8
$ws.Range("M5").Value = -181.8

$ws.Range("H31").Value = 6083.1
$ws.Range("I31").Value = 6083.1
$ws.Range("K31").Value = 6083.1
$ws.Range("M31").Value = -5789.1

$ws.Range("H32").Value = 3709457.5
$ws.Range("I32").Value = 4548702.5
$ws.Range("J32").Value = 16779.2
$ws.Range("K32").Value = 4548702.5
$ws.Range("L32").Value = 16779.2
$ws.Range("M32").Value = -4548415.5
$ws.Range("N32").Value = -17353.2

$ws.Range("H63").Value = 3767.1667
$ws.Range("J63").Value = 2000
$ws.Range("L63").Value = 2000
$ws.Range("N63").Value = -3372

$ws.Range("H66").Value = 3767.1667
$ws.Range("J66").Value = 2000
$ws.Range("L66").Value = 10000
$ws.Range("N66").Value = -16864

$ws.Range("H102").Value = 2725
$ws.Range("I102").Value = 2725
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2725
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1103
$ws.Range("N102").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 314.83334
$ws.Range("I4").Value = 293.8
$ws.Range("K4").Value = 293.8
$ws.Range("M4").Value = -178.8

$ws.Range("H20").Value = 26608.229
$ws.Range("J20").Value = 39510.156
$ws.Range("L20").Value = 39510.156
$ws.Range("N20").Value = -40004.156

$ws.Range("H105").Value = 1268.5416
$ws.Range("I105").Value = 1046.5883
$ws.Range("K105").Value = 1046.5883
$ws.Range("M105").Value = 700.4117000000001

$ws.Range("H138").Value = 65000
$ws.Range("I138").Value = 65000
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 65000
$ws.Range("L138").Value = 0
$ws.Range("M138").Value = -59860
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 10000
$ws.Range("I4").Value = 10000
$ws.Range("K4").Value = 10000
$ws.Range("M4").Value = -9888

$ws.Range("H31").Value = 148344.14
$ws.Range("I31").Value = 251204.75
$ws.Range("J31").Value = 30789.143
$ws.Range("K31").Value = 251204.75
$ws.Range("L31").Value = 30789.143
$ws.Range("M31").Value = -250909.75
$ws.Range("N31").Value = -31379.143

$ws.Range("H34").Value = 148344.14
$ws.Range("I34").Value = 251204.75
$ws.Range("J34").Value = 30789.143
$ws.Range("K34").Value = 251204.75
$ws.Range("L34").Value = 30789.143
$ws.Range("M34").Value = -251002.75
$ws.Range("N34").Value = -31193.143

$ws.Range("H62").Value = 7997.5
$ws.Range("J62").Value = 7997.5
$ws.Range("L62").Value = 7997.5
$ws.Range("N62").Value = -9245.5

$ws.Range("H65").Value = 7997.5
$ws.Range("J65").Value = 7997.5
$ws.Range("L65").Value = 39987.5
$ws.Range("N65").Value = -46227.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 4340
$ws.Range("J113").Value = 1093
$ws.Range("L113").Value = 3279
$ws.Range("N113").Value = -7619

$ws.Range("H122").Value = 10248567
$ws.Range("J122").Value = 2361533.5
$ws.Range("L122").Value = 21253801.5
$ws.Range("N122").Value = -21258701.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 73.61905
$ws.Range("I2").Value = 84.5
$ws.Range("J2").Value = 38.8
$ws.Range("K2").Value = 84.5
$ws.Range("L2").Value = 38.8
$ws.Range("M2").Value = 28.5
$ws.Range("N2").Value = -264.8

$ws.Range("H113").Value = 3505.44
$ws.Range("I113").Value = 2113.8667
$ws.Range("K113").Value = 2113.8667
$ws.Range("M113").Value = 56.13329999999996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 5326.5
$ws.Range("J22").Value = 5326.5
$ws.Range("L22").Value = 5326.5
$ws.Range("N22").Value = -5916.5

$ws.Range("H27").Value = 5326.5
$ws.Range("J27").Value = 5326.5
$ws.Range("L27").Value = 5326.5
$ws.Range("N27").Value = -5540.5

$ws.Range("H46").Value = 3417.3635
$ws.Range("J46").Value = 3417.3635
$ws.Range("L46").Value = 3417.3635
$ws.Range("N46").Value = -3793.3635

$ws.Range("H55").Value = 1721.5
$ws.Range("I55").Value = 1537.1818
$ws.Range("J55").Value = 1946.7778
$ws.Range("K55").Value = 1537.1818
$ws.Range("L55").Value = 1946.7778
$ws.Range("M55").Value = -1364.1818
$ws.Range("N55").Value = -2292.7778

$ws.Range("H68").Value = 14229.75
$ws.Range("I68").Value = 18306.334
$ws.Range("K68").Value = 18306.334
$ws.Range("M68").Value = -17557.334

$ws.Range("H71").Value = 14229.75
$ws.Range("I71").Value = 18306.334
$ws.Range("K71").Value = 91531.67
$ws.Range("M71").Value = -87787.67

$ws.Range("H82").Value = 2349.6538
$ws.Range("I82").Value = 2644.4375
$ws.Range("J82").Value = 1878
$ws.Range("K82").Value = 2644.4375
$ws.Range("L82").Value = 1878
$ws.Range("M82").Value = -2283.4375
$ws.Range("N82").Value = -2600

$ws.Range("H85").Value = 2349.6538
$ws.Range("I85").Value = 2644.4375
$ws.Range("J85").Value = 1878
$ws.Range("K85").Value = 2644.4375
$ws.Range("L85").Value = 1878
$ws.Range("M85").Value = -1396.4375
$ws.Range("N85").Value = -4374

$ws.Range("H100").Value = 3168.1667
$ws.Range("I100").Value = 2857.2
$ws.Range("K100").Value = 2857.2
$ws.Range("M100").Value = -2316.2

$ws.Range("H122").Value = 5899.9062
$ws.Range("I122").Value = 5915.95
$ws.Range("J122").Value = 5873.1665
$ws.Range("K122").Value = 17747.85
$ws.Range("L122").Value = 17619.4995
$ws.Range("M122").Value = -15297.85
$ws.Range("N122").Value = -22519.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1750.5
$ws.Range("J96").Value = 1800.6
$ws.Range("L96").Value = 1800.6
$ws.Range("N96").Value = -4546.6
